# Insert a new data row into the "Hortaliza, Feria Lagunitas de Puerto Montt -
# Pepino ensalada" sheet. The diff shows a brand-new record (date 2023-05-29,
# i.e. serial 45075) inserted right after the existing row 423, which pushes
# every following row (old 424-439) down by one (new 425-440) and grows the
# used range from A1:R439 to A1:R440.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 424-439 down to 425-440, leaving a fresh blank row 424.
$ws.Rows.Item(424).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(424, 1).Value = 4
$ws.Cells.Item(424, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(424, 3).Value = "Los Lagos"
$ws.Cells.Item(424, 4).Value = 45075
$ws.Cells.Item(424, 5).Value = 10
$ws.Cells.Item(424, 6).Value = 100112043
$ws.Cells.Item(424, 7).Value = "Pepino ensalada"
$ws.Cells.Item(424, 8).Value = "Sin especificar"
$ws.Cells.Item(424, 9).Value = "Primera"
$ws.Cells.Item(424, 10).Value = 120
$ws.Cells.Item(424, 11).Value = 16000
$ws.Cells.Item(424, 12).Value = 16000
$ws.Cells.Item(424, 13).Value = 16000
$ws.Cells.Item(424, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(424, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(424, 16).Value = 267
$ws.Cells.Item(424, 17).Value = 60
$ws.Cells.Item(424, 18).Value = "Hortaliza"
